# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Sat Apr 13 03:06:42 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.527.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -5.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.212.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -8.78%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -5.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.47"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -13.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.203.81"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -8.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -11.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -13.50%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -11.63%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -15.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.09"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -17.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000241"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -12.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.723.74"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -9.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.440.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.208.54"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -8.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "534.20"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -13.19%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.09"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -16.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.96"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -15.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.752"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -14.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.65"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -14.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -13.23%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -16.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.12"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -17.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -18.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.94"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -13.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.95"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -14.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.53"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -16.75%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -13.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "530.46"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -13.86%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -20.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.59"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -18.26%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.88"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0418"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -12.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0848"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -15.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -17.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.125"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -13.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.889.86"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -14.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.59"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -26.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0582"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -20.71%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -17.50%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -21.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.66"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -20.41%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -18.90%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.28"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.76%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.112"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -14.10%  "
